$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column widths: split the old "B:G" block so column D gets its own wider width ---
$ws.Range("D1").ColumnWidth = 34.5

# --- Move the "TimeSlice" header (with its formatting) from J3 down to H5 ---
$ws.Range("J3").Copy()
$ws.Range("H5").PasteSpecial(-4122)
$ws.Range("H5").Value = "TimeSlice"

# Clear out the now-empty J3 cell (value + formatting) so row 3 collapses away
$ws.Range("J3").Clear()
$ws.Rows(3).AutoFit()

# --- D6: replace the single year 2030 with the full list of years ---
$ws.Range("D6").Value = "2025,2030,2035,2040,2045,2050"

# --- Update the active selection to D8 (matches the saved selection in the file) ---
$ws.Range("D8").Select() | Out-Null
